$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the requirement-related column headers in row 1 to be more
# explicit (e.g. "Reference" -> "Requirement Reference"). The underlying
# Excel table ("Tabelle1") automatically picks up the new header text for
# its column names.
$ws.Range("A1").Value = "Requirement Reference"
$ws.Range("B1").Value = "Requirement Summary"
$ws.Range("C1").Value = "Requirement Description"
$ws.Range("E1").Value = "Requirement Compliance Status"
$ws.Range("F1").Value = "Requirement Compliance Comment"
$ws.Range("G1").Value = "Requirement Completion Progress"

# The new header text for column G is noticeably longer, so widen the
# column to fit it (matches Excel's own best-fit behaviour).
$ws.Columns("G").AutoFit() | Out-Null

# Leave the selection on G2, matching where the cursor ended up.
$ws.Range("G2").Select() | Out-Null
